$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Modelo" header in F1, matching the style of the existing headers (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update numeric values in row 2 (MSE, R2, MAE)
$ws.Range("B2").Value = 0.09406460280887817
$ws.Range("C2").Value = 0.9994396839998366
$ws.Range("D2").Value = 0.2170285180861791

# Add the new model description in F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5, n_estimators=50))])"
